# Week 5 work: drop the "Recycling Process" write-up (and its now-unused
# list numbering/style) from the end of the flowchart document; the
# diagram/header/footer relationship ids shift down to fill the gap left
# by the removed numbering.xml relationship.

$d = $word.ActiveDocument

# Pull the whole package (flat OPC) so we can edit document.xml,
# word/_rels/document.xml.rels, numbering.xml and styles.xml together -
# Word keeps their relationship ids/content in lock-step.
$xml = $d.Content.WordOpenXML

# ---------------------------------------------------------------------
# 1. Remove the "Recycling Process:" paragraph and the six bullet items
#    that follow it, right before the closing </w:body> sectPr. The
#    empty paragraph just before stays untouched.
# ---------------------------------------------------------------------
$blockStart = $xml.IndexOf('<w:p w14:paraId="277206F2"')
if ($blockStart -lt 0) {
    throw "Could not find the start of the Recycling Process block"
}
$blockEnd = $xml.IndexOf('<w:sectPr')
if ($blockEnd -lt 0) {
    throw "Could not find sectPr after the Recycling Process block"
}
$xml = $xml.Substring(0, $blockStart) + $xml.Substring($blockEnd)

# ---------------------------------------------------------------------
# 2. Drop the now-orphaned numbering.xml part and its relationship.
# ---------------------------------------------------------------------
$relStart = $xml.IndexOf('<pkg:part pkg:name="/word/_rels/document.xml.rels"')
$relDataStart = $xml.IndexOf('<pkg:xmlData>', $relStart)
$relDataEnd = $xml.IndexOf('</pkg:xmlData>', $relDataStart)
if ($relStart -lt 0 -or $relDataStart -lt 0 -or $relDataEnd -lt 0) {
    throw "Could not locate document.xml.rels part"
}
$relsXml = $xml.Substring($relDataStart + 13, $relDataEnd - ($relDataStart + 13))
$numRelPattern = '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/numbering" Target="numbering.xml"/>'
if ($relsXml.IndexOf($numRelPattern) -lt 0) {
    throw "Could not find the numbering relationship entry"
}
$relsXml = $relsXml.Replace($numRelPattern, "")
$xml = $xml.Substring(0, $relDataStart + 13) + $relsXml + $xml.Substring($relDataEnd)

$numPartStart = $xml.IndexOf('<pkg:part pkg:name="/word/numbering.xml"')
if ($numPartStart -lt 0) {
    throw "Could not find numbering.xml part"
}
$numPartEnd = $xml.IndexOf('</pkg:part>', $numPartStart)
if ($numPartEnd -lt 0) {
    throw "Could not find end of numbering.xml part"
}
$numPartEnd = $numPartEnd + ('</pkg:part>').Length
$xml = $xml.Substring(0, $numPartStart) + $xml.Substring($numPartEnd)

# ---------------------------------------------------------------------
# 3. Remove the (now unused) ListParagraph style definition.
# ---------------------------------------------------------------------
$styleStart = $xml.IndexOf('<w:style w:type="paragraph" w:styleId="ListParagraph">')
if ($styleStart -lt 0) {
    throw "Could not find the ListParagraph style"
}
$styleEnd = $xml.IndexOf('</w:style>', $styleStart)
if ($styleEnd -lt 0) {
    throw "Could not find end of the ListParagraph style"
}
$styleEnd = $styleEnd + ('</w:style>').Length
$xml = $xml.Substring(0, $styleStart) + $xml.Substring($styleEnd)

# ---------------------------------------------------------------------
# 4. With rId1 (numbering) gone, every remaining rIdN shifts down to
#    rId(N-1) - Word renumbers relationship ids to stay contiguous.
# ---------------------------------------------------------------------
$matches = [regex]::Matches($xml, 'rId(\d+)')
$renumbered = ""
$lastEnd = 0
foreach ($m in $matches) {
    $n = [int]$m.Groups[1].Value
    $renumbered += $xml.Substring($lastEnd, $m.Index - $lastEnd)
    $renumbered += ("rId" + ($n - 1))
    $lastEnd = $m.Index + $m.Length
}
$renumbered += $xml.Substring($lastEnd)
$xml = $renumbered

# Write the edited package back - InsertXML on the whole-document range
# replaces the entire content (and its associated parts), same as
# pasting a pkg:package fragment in real Word automation.
$d.Content.InsertXML($xml)
